# Rebuild rows 2-42 of the "memory / living_rooms" trial table.
# The task now contains 41 trials (was 28): more trial rows were added
# and the mix of old/new stimuli was rebalanced toward more "new" items,
# per the commit message ("increased the proportion of new images in
# memory task") and elaborated sanity-check (cond_cat) tagging.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A subject_id, B task, C block_total, D block_scene, E trial_block,
# F trial_total, G target_cat, H category, I cond_cat, J cond_mem,
# K correct_answer, L stimulus, M conceptual, N perceptual, O typicality,
# P n, Q p_typicality, R p_conceptual, S p_perceptual
$rows = @(
    @(18, 'memory', 3, 1, 1, 161, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_7lz7m.png', 51.5531914893617, 32.87234042553192, 42.21276595744681, 47, 3, 3, 3),
    @(18, 'memory', 3, 1, 2, 162, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_kost0.png', 63.09090909090909, 42.77272727272727, 52.93181818181819, 44, 5, 5, 5),
    @(18, 'memory', 3, 1, 3, 163, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_bbs77.png', 31.64444444444445, 21.26666666666667, 26.45555555555556, 45, 2, 2, 2),
    @(18, 'memory', 3, 1, 4, 164, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_b17ma.png', 23.0625, 13.375, 18.21875, 48, 1, 1, 1),
    @(18, 'memory', 3, 1, 5, 165, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_wgkqa.png', 87.25581395348837, 71.13953488372093, 79.19767441860465, 43, 10, 10, 10),
    @(18, 'memory', 3, 1, 6, 166, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_6zz63.png', 87.66666666666667, 70.6, 79.13333333333333, 45, 9, 10, 10),
    @(18, 'memory', 3, 1, 7, 167, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_di6f0.png', 94.04347826086956, 83.34782608695652, 88.69565217391303, 46, 10, 10, 10),
    @(18, 'memory', 3, 1, 8, 168, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_f63yi.png', 85.275, 68.475, 76.875, 40, 9, 9, 9),
    @(18, 'memory', 3, 1, 9, 169, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_zh8ms.png', 59.82608695652174, 39.43478260869565, 49.6304347826087, 46, 4, 4, 4),
    @(18, 'memory', 3, 1, 10, 170, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_6a0hu.png', 61.275, 42.025, 51.65, 40, 4, 4, 4),
    @(18, 'memory', 3, 1, 11, 171, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_4o8l0.png', 46.02173913043478, 31.45652173913043, 38.73913043478261, 46, 3, 3, 3),
    @(18, 'memory', 3, 1, 12, 172, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_16kib.png', 80.97727272727273, 61.11363636363637, 71.04545454545455, 44, 8, 8, 8),
    @(18, 'memory', 3, 1, 13, 173, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_15bss.png', 88.42222222222222, 75.35555555555555, 81.88888888888889, 45, 10, 10, 10),
    @(18, 'memory', 3, 1, 14, 174, 'living_rooms', $null, $null, 'catch', 'f', 'stimuli/catch_01.jpg', $null, $null, $null, $null, $null, $null, $null),
    @(18, 'memory', 3, 1, 15, 175, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_bj99b.png', 82.79069767441861, 65.46511627906976, 74.12790697674419, 43, 8, 8, 8),
    @(18, 'memory', 3, 1, 16, 176, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_3jnt7.png', 49.52272727272727, 35.25, 42.38636363636364, 44, 3, 3, 3),
    @(18, 'memory', 3, 1, 17, 177, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_fmgjx.png', 79.9, 56.975, 68.4375, 40, 7, 7, 7),
    @(18, 'memory', 3, 1, 18, 178, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_lgxzn.png', 73.11363636363636, 49.97727272727273, 61.54545454545455, 44, 6, 6, 6),
    @(18, 'memory', 3, 1, 19, 179, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_pey7u.png', 30.34883720930232, 20.34883720930232, 25.34883720930232, 43, 1, 2, 2),
    @(18, 'memory', 3, 1, 20, 180, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_5mw7y.png', 72.6590909090909, 50.86363636363637, 61.76136363636364, 44, 6, 6, 6),
    @(18, 'memory', 3, 1, 21, 181, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_0jzz7.png', 84.85106382978724, 68.87234042553192, 76.86170212765958, 47, 9, 9, 9),
    @(18, 'memory', 3, 1, 22, 182, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_wz6x5.png', 68.3695652173913, 48.47826086956522, 58.42391304347826, 46, 5, 5, 5),
    @(18, 'memory', 3, 1, 23, 183, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_qrc78.png', 76.2, 59.875, 68.0375, 40, 7, 7, 7),
    @(18, 'memory', 3, 1, 24, 184, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_s2zoe.png', 64.71428571428571, 44.90476190476191, 54.80952380952381, 42, 5, 5, 5),
    @(18, 'memory', 3, 1, 25, 185, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_eiu3c.png', 65.1590909090909, 46.22727272727273, 55.69318181818181, 44, 5, 5, 5),
    @(18, 'memory', 3, 1, 26, 186, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_w8yhd.png', 55.74418604651163, 38.90697674418605, 47.32558139534883, 43, 4, 4, 4),
    @(18, 'memory', 3, 1, 27, 187, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_lpas9.png', 59.36585365853659, 39.09756097560975, 49.23170731707317, 41, 4, 4, 4),
    @(18, 'memory', 3, 1, 28, 188, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_xu1p3.png', 75.27659574468085, 56.68085106382978, 65.97872340425532, 47, 7, 7, 7),
    @(18, 'memory', 3, 1, 29, 189, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_xbtev.png', 13.68181818181818, 8.568181818181818, 11.125, 44, 1, 1, 1),
    @(18, 'memory', 3, 1, 30, 190, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_0kqc0.png', 43.74468085106383, 27.14893617021277, 35.4468085106383, 47, 2, 2, 2),
    @(18, 'memory', 3, 1, 31, 191, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_eh0no.png', 53.66666666666666, 36.02564102564103, 44.84615384615385, 39, 3, 3, 3),
    @(18, 'memory', 3, 1, 32, 192, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_il020.png', 18.85416666666667, 16.16666666666667, 17.51041666666667, 48, 1, 1, 1),
    @(18, 'memory', 3, 1, 33, 193, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_cehin.png', 78.86363636363636, 60.02272727272727, 69.44318181818181, 44, 7, 7, 7),
    @(18, 'memory', 3, 1, 34, 194, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_89dvt.png', 81.09756097560975, 64.6829268292683, 72.89024390243902, 41, 8, 8, 8),
    @(18, 'memory', 3, 1, 35, 195, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_pjfx6.png', 32.23404255319149, 26.59574468085106, 29.41489361702127, 47, 2, 2, 2),
    @(18, 'memory', 3, 1, 36, 196, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_314bq.png', 37.08888888888889, 20.04444444444444, 28.56666666666667, 45, 2, 2, 2),
    @(18, 'memory', 3, 1, 37, 197, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_emh91.png', 82.06666666666666, 63.33333333333334, 72.7, 45, 8, 8, 8),
    @(18, 'memory', 3, 1, 38, 198, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_abobq.png', 75.1842105263158, 54.13157894736842, 64.65789473684211, 38, 6, 6, 6),
    @(18, 'memory', 3, 1, 39, 199, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_xy930.png', 70.5952380952381, 49.47619047619047, 60.03571428571429, 42, 6, 6, 6),
    @(18, 'memory', 3, 1, 40, 200, 'living_rooms', 'living_rooms', $null, 'new', 'f', 'stimuli/img_dg5h7.png', 88.72093023255815, 76.06976744186046, 82.3953488372093, 43, 10, 10, 10),
    @(18, 'memory', 3, 1, 41, 201, 'living_rooms', 'living_rooms', 'target', 'old', 'j', 'stimuli/img_amsgw.png', 86.08510638297872, 65.95744680851064, 76.02127659574468, 47, 9, 9, 9)
)

$nRows = $rows.Count
$nCols = 19
$data = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$lastRow = 1 + $nRows
$ws.Range("A2:S$lastRow").Value = $data

Write-Output ("New used range: " + $ws.UsedRange.Address())
